# The deck ships with two embedded theme parts: the design that is
# actually bound to the (single) slide master - "Integral" / "Red
# Violet" - and a spare, unused "Office Theme" part that simply rides
# along in the package. Switching the presentation's Design from
# "Integral" back to the built-in "Office Theme" (Design gallery)
# re-colours the live theme with the stock Office palette.
#
# RGB() isn't available in this host, so colours are written as the
# packed 0xBBGGRR integer PowerPoint's ColorFormat.RGB expects
# (r + g*256 + b*65536).

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Target palette: the built-in "Office" colour scheme.
#  1 Dark1   2 Light1  3 Dark2   4 Light2
#  5 Accent1 6 Accent2 7 Accent3 8 Accent4 9 Accent5 10 Accent6
# 11 Hyperlink 12 FollowedHyperlink
$colorScheme.Colors(1).RGB  = 0x000000   # dk1
$colorScheme.Colors(2).RGB  = 0xFFFFFF   # lt1
$colorScheme.Colors(3).RGB  = 0x6A5444   # dk2 (44546A, BGR-packed)
$colorScheme.Colors(4).RGB  = 0xE6E6E7   # lt2 (E7E6E6)
$colorScheme.Colors(5).RGB  = 0xD59B5B   # accent1 (5B9BD5)
$colorScheme.Colors(6).RGB  = 0x317DED   # accent2 (ED7D31)
$colorScheme.Colors(7).RGB  = 0xA5A5A5   # accent3 (A5A5A5)
$colorScheme.Colors(8).RGB  = 0x00C0FF   # accent4 (FFC000)
$colorScheme.Colors(9).RGB  = 0xC47244   # accent5 (4472C4)
$colorScheme.Colors(10).RGB = 0x47AD70   # accent6 (70AD47)
$colorScheme.Colors(11).RGB = 0xC16305   # hlink (0563C1)
$colorScheme.Colors(12).RGB = 0x724F95   # folHlink (954F72)
